## Applies "Added more information to timeline": appends 10 more rows
## (rows 6-15) to the Job Applications timeline sheet, alternating between
## two job entries (Software Developer Co-Op @ TD / Web Developer Intern @
## Pathcore), each with its own Skills blurb and hyperlinked Url cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$linkedInUrl = "https://www.linkedin.com/jobs/view/1865907079/?eBP=JYMBII_JOBS_HOME_ORGANIC&recommendedFlavor=SCHOOL_RECRUIT&refId=f6e15b88-0061-4368-b3ce-fe439cee172c&trk=d_flagship3_jobs_discovery_jymbii"
$indeedUrl   = "https://ca.indeed.com/viewjob?jk=45e7adfb4d34664e&tk=1e8k9749r0gc1000&from=serp&vjs=3"

$skillsA = "HTML CSS Javascript Java C++ Objective c GIT Swift "
$skillsB = "Python HTML CSS Javascript Java GIT "

for ($row = 6; $row -le 15; $row++) {

    $ws.Range("A" + $row).Value = $row - 1

    # "Date" column is text that looks like an ISO date ("2020-05-21"); force
    # the cell to Text first so Excel doesn't silently coerce it to a date
    # serial number the way it would for a bare Value assignment.
    $ws.Range("B" + $row).NumberFormat = "@"
    $ws.Range("B" + $row).Value = "2020-05-21"

    if (($row % 2) -eq 0) {
        $ws.Range("C" + $row).Value = "Software Developer Co-Op - TD Innovation Lab"
        $ws.Range("D" + $row).Value = "TD"
        $ws.Range("E" + $row).Value = $skillsA
        $ws.Range("F" + $row).Value = $linkedInUrl
        $ws.Hyperlinks.Add($ws.Range("F" + $row), $linkedInUrl)
    } else {
        $ws.Range("C" + $row).Value = "Web Developer Intern"
        $ws.Range("D" + $row).Value = "Pathcore"
        $ws.Range("E" + $row).Value = $skillsB
        $ws.Range("F" + $row).Value = $indeedUrl
        $ws.Hyperlinks.Add($ws.Range("F" + $row), $indeedUrl)
    }
}

# Re-apply the formatting (and only the formatting) from the first data row
# onto the newly-added rows, so the new cells line up with the existing
# borderless/bold/centered "index" style (col A) and the Hyperlink style
# (col F) rather than whatever incidental formatting the writes above left
# behind.
$ws.Range("A2").Copy()
$ws.Range("A6:A15").PasteSpecial(-4122)

$ws.Range("B2").Copy()
$ws.Range("B6").PasteSpecial(-4122)
$ws.Range("B7").PasteSpecial(-4122)
$ws.Range("B8").PasteSpecial(-4122)
$ws.Range("B9").PasteSpecial(-4122)
$ws.Range("B10").PasteSpecial(-4122)
$ws.Range("B11").PasteSpecial(-4122)
$ws.Range("B12").PasteSpecial(-4122)
$ws.Range("B13").PasteSpecial(-4122)
$ws.Range("B14").PasteSpecial(-4122)
$ws.Range("B15").PasteSpecial(-4122)

$ws.Range("F2").Copy()
$ws.Range("F6:F15").PasteSpecial(-4122)

$excel.CutCopyMode = $false
